# Update the "department" column (C) with more granular department names,
# replacing the single generic "FACULTY OF BUSINESS & TECHNOLOGY" value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Business department courses
$ws.Range("C2:C9").Value = "Business"

# Information Technology department courses
$ws.Range("C10:C12").Value = "Information Technology"

# Building and Construction department course
$ws.Range("C13").Value = "Building and Construction"

# Packaged courses
$ws.Range("C14:C19").Value = "Packages"
